$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A9").Value = "Natalia Pelaez Arboleda"
$ws.Range("B9").Value = "npelaez@bancolombia.com.co"
$ws.Range("C9").Value = "BPM07 - FLASH"
$ws.Range("D9").Value = "AW1176001_BIZAGIECUCOL_TEST"

$ws.Range("A9:D9").Select()
